# ASA Logical Device added
# This script reshapes the "# Logical Device Configuration" block on the
# "FXOS DC2 Settings" sheet (sheet index 2) into a full ASA + FTD logical
# device configuration block, and renames "Portchannel" -> "Port-channel"
# in the Port-Channel Configuration table above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 1. Port-Channel Configuration table: rename label -------------------
# Rows 57-59 column A used to read "Portchannel"; now reads "Port-channel".
$ws.Range("A57").Value = "Port-channel"
$ws.Range("A58").Value = "Port-channel"
$ws.Range("A59").Value = "Port-channel"

# Row 57 used to list every member interface (H:O); the last three
# (Ethernet1/6, Ethernet1/7, Ethernet1/8) are no longer listed there.
$ws.Range("M57:O57").ClearContents()

# --- 2. Logical Device Configuration: headers -----------------------------
# Shared header layout reused for both the ASA block (row 61) and the new
# FTD block (row 64).
$headerCols = @{
    2  = "Slot Number"
    3  = "Hostname"
    4  = "Software Version"
    5  = "Firewall Mode"
    6  = "Admin Password"
    7  = "Management IP"
    8  = "Subnet Mask"
    9  = "Default Gateway"
    10 = "Management Interface"
    11 = "Nameif"
    12 = "Description"
    13 = "Interface 1"
    14 = "Nameif"
    15 = "Description"
    16 = "Interface 2"
    17 = "Nameif"
    18 = "Description"
    19 = "Interface 3"
    20 = "Nameif"
    21 = "Description"
    22 = "Interface 4"
    23 = "Nameif"
    24 = "Description"
}

# Row 61: ASA section header
$ws.Cells.Item(61, 1).Value = "# ASA Logical Device Configuration"
$ws.Cells.Item(61, 1).Font.Bold = $true
foreach ($col in $headerCols.Keys) {
    $cell = $ws.Cells.Item(61, $col)
    $cell.Value = $headerCols[$col]
    $cell.Font.Bold = $true
}

# Row 62: ASA data row
$asaRow = @{
    1  = "ASA"
    2  = 1
    3  = "ASA1"
    4  = "9.12.1"
    5  = "routed"
    6  = "cisco123"
    7  = "192.168.1.10"
    8  = "255.255.255.0"
    9  = "192.168.1.1"
    10 = "Port-channel 30"
    11 = "management"
    12 = "management interface"
    13 = "Port-channel 10"
    14 = "inside"
    15 = "inside interface"
    16 = "Port-channel 20"
    17 = "outside"
    18 = "outside interface"
    19 = "Ethernet2/7"
    20 = "DMZ"
    21 = "DMZ interface"
    22 = "Ethernet2/8"
    23 = "DMZ2"
    24 = "DMZ2 interface"
}
foreach ($col in $asaRow.Keys) {
    $ws.Cells.Item(62, $col).Value = $asaRow[$col]
}

# Row 64: FTD section header (same layout as row 61)
$ws.Cells.Item(64, 1).Value = "# FTD Logical Device Configuration"
$ws.Cells.Item(64, 1).Font.Bold = $true
foreach ($col in $headerCols.Keys) {
    $cell = $ws.Cells.Item(64, $col)
    $cell.Value = $headerCols[$col]
    $cell.Font.Bold = $true
}

# Row 65: FTD data row (re-uses existing Ethernet1/6-1/8 strings)
$ftdRow = @{
    1  = "FTD"
    2  = 2
    3  = "FTD1"
    4  = "6.2.3.83"
    5  = "routed"
    6  = "cisco123"
    7  = "192.168.1.20"
    8  = "255.255.255.0"
    9  = "192.168.1.1"
    10 = "Ethernet1/8"
    11 = "management"
    12 = "management interface"
    13 = "Ethernet1/6"
    14 = "inside"
    15 = "inside interface"
    16 = "Ethernet1/7"
    17 = "outside"
    18 = "outside interface"
}
foreach ($col in $ftdRow.Keys) {
    $ws.Cells.Item(65, $col).Value = $ftdRow[$col]
}

# --- 3. Column widths -------------------------------------------------------
# The wider ASA/FTD tables push a handful of columns to re-autofit; reproduce
# the resulting "best fit" widths (expressed here as the COM character-width
# value that this engine's 5/6-char padding model round-trips to the target
# stored width).
$ws.Columns.Item(9).ColumnWidth  = 20.877604166666668  # col I  -> 21.7109375
$ws.Columns.Item(12).ColumnWidth = 20.736979166666668  # col L  -> 21.5703125
$ws.Columns.Item(16).ColumnWidth = 12.736979166666666  # col P  -> 13.5703125
$ws.Columns.Item(17).ColumnWidth = 11.307291666666666  # col Q  -> 12.140625
$ws.Columns.Item(18).ColumnWidth = 15.451822916666666  # col R  -> 16.28515625
$ws.Columns.Item(19).ColumnWidth = 13.736979166666666  # col S  -> 14.5703125

# --- 4. View state ---------------------------------------------------------
$ws.Activate()
$ws.Range("H75").Select()
